# Asservissement - doc to establish transfert function
# This script edits the "acc" sheet to move from a 20-bit model to a
# 19-bit model (one fewer bit), swaps the neg_max/neg_min labels that were
# previously inverted, rescales several sample inputs, and removes the
# last two helper rows (the old bit-20 row and the trailing g/-8 row,
# which is folded up into row 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("acc")

# --- Swap the neg_max / neg_min labels in G7 / G8 (they were reversed) ---
$ws.Range("G7").Value = "neg_min"
$ws.Range("G8").Value = "neg_max"

# --- Shrink the summation ranges now that the model has one fewer bit ---
$ws.Range("H4").Formula = "=SUM(D4:D22)"
$ws.Range("H7").Formula = "=SUM(D4:D23)"
$ws.Range("D28").Formula = "=SUM(D4:D23)"
$ws.Range("D30").Formula = "=SUM(D4:D22)"

# --- Remove the old "bit 20" row (B24/D24); row stays, just those 2 cells go ---
$ws.Range("B24:D24").ClearContents()

# --- G20 becomes a plain pasted value instead of a live SUM formula ---
$ws.Range("G20").Formula = "524287"

# --- Rescaled sample inputs in column G ---
$ws.Range("G21").Value = 400000
$ws.Range("G22").Value = 300000
$ws.Range("G23").Value = 250000
$ws.Range("G28").Value = 524288
$ws.Range("G29").Value = 600000
$ws.Range("G30").Value = 700000
$ws.Range("G31").Value = 750000
$ws.Range("G32").Value = 800000
$ws.Range("G33").Value = 900000
$ws.Range("G34").Value = 1048575

# --- Fold the old row 35 (I35/J35, the "g"/-8 helper) up into row 34 ---
$ws.Range("I34").Formula = "=-8"
$ws.Range("J34").Value = "g"
$ws.Range("J34").HorizontalAlignment = -4131

# Now that row 34 carries the final content, delete the now-redundant row 35
# (shifts nothing else up, since row 35 was the last used row).
$ws.Rows("35:35").EntireRow.Delete()

# --- Update the cursor / selection like the saved file shows ---
$ws.Range("I38").Select()

$wb.Save()
